# Revert "Merge pull request #5 ... debug-compute-sensitivity-code"
# i.e. undo the changes that PR introduced, restoring the prior state of
# Install/Simulation_Output_Settings.xlsx.

$wb = $excel.ActiveWorkbook

$wsTotal  = $wb.Worksheets.Item("TotalEnergy")
$wsMeters = $wb.Worksheets.Item("Meters")
$wsSources = $wb.Worksheets.Item("Sources")

# --- TotalEnergy: drop the "electricityCooling" row (was A2) -------------
# A1=Total Energy Output Variable(24), A2=electricityCooling(25, removed),
# A3=Total Site Energy[GJ](27)->becomes A2, A4=Electricity Total End Uses[GJ](31)->becomes A3
$wsTotal.Rows.Item(2).Delete()

# --- Sources: drop the carbon-equivalent / fans / pumps / water-system ----
# helper rows that the debug-sensitivity branch had added in column D
# (D9:D14 and D16), restoring the shorter reference table.
$wsSources.Range("D9:D14").ClearContents()
$wsSources.Range("D16").ClearContents()

# --- column width tweaks (best-effort match of the reverted layout) -------
$wsTotal.Columns.Item(2).ColumnWidth = 23.333333333333332
$wsTotal.Columns.Item(3).ColumnWidth = 18.5
$wsTotal.Columns.Item(4).ColumnWidth = 32.5

$wsMeters.Columns.Item(1).ColumnWidth = 28.5
$wsMeters.Columns.Item(2).ColumnWidth = 21.0

$wsSources.Columns.Item(1).ColumnWidth = 40.5

# --- restore the prior active sheet / selections ---------------------------
# Before the merge, Sources!E11 and TotalEnergy!A3 were the last selections
# on those (now inactive) tabs, and Meters!A3 is the active selection/tab
# (so it must be the LAST sheet activated, to end up as tabSelected).
$wsTotal.Activate()
$wsTotal.Range("A3").Select()

$wsSources.Activate()
$wsSources.Range("E11").Select()

$wsMeters.Activate()
$wsMeters.Range("A3").Select()

Write-Output "edit applied"
